# Updated legacy GSC export data:
#  - drop the oldest date row (2025-10-02) and shift all remaining
#    rows up by one
#  - append a new row for the next day (2025-12-31) with zero counts
#
# This mirrors the "Chart" worksheet, which is the one holding the
# per-day HTTPS/Non-HTTPS URL counts (columns A:C, rows 2:91).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Find the last populated data row (row 1 is the header "Date" / "Non-HTTPS URLs" / "HTTPS URLs").
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp
$newRow  = $lastRow  # the shifted data now ends one row earlier; this row becomes the new entry

# 1) Shift rows 3..91 (A:C) up into rows 2..90, which drops the old
#    first data row (2025-10-02, C=67) and moves every other row up
#    by one, exactly like removing the earliest date from the export.
$srcRange = $ws.Range("A3:C" + $lastRow)
$dstRange = $ws.Range("A2:C" + ($lastRow - 1))
$srcRange.Copy($dstRange)

# 2) Build the new date string "2025-12-31" through a formula so Excel
#    does not auto-convert the literal text into a date serial number,
#    then paste only the computed value so the cell stays a plain text
#    cell (same as every other date cell) without picking up a new
#    number format/style.
$scratch = $ws.Range("Z1")
$scratch.Formula = '=TEXT(DATE(2025,12,31),"yyyy-mm-dd")'
$scratch.Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4163)  # xlPasteValues
$scratch.ClearContents()

# 3) New row's counts: no data yet for the new day.
$ws.Range("B" + $newRow).Value = 0
$ws.Range("C" + $newRow).Value = 0
